$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 11).Value = 0.0652635207166864
$ws.Cells.Item(3, 2).Value = 7.869792681105293
$ws.Cells.Item(3, 3).Value = -10.45921331889471
$ws.Cells.Item(3, 4).Value = -2.319131318894706
$ws.Cells.Item(3, 5).Value = -0.5308223188947059
$ws.Cells.Item(3, 6).Value = -3.708941318894706
$ws.Cells.Item(3, 7).Value = -4.004270318894706
$ws.Cells.Item(3, 8).Value = -1.489825318894706
$ws.Cells.Item(3, 9).Value = -2.419547318894706
$ws.Cells.Item(3, 10).Value = -2.184539318894706
$ws.Cells.Item(3, 11).Value = -2.491675318894706
$ws.Cells.Item(4, 2).Value = -17.928363649752
$ws.Cells.Item(4, 3).Value = -9.788281649751999
$ws.Cells.Item(4, 4).Value = -7.999972649751999
$ws.Cells.Item(4, 5).Value = -11.178091649752
$ws.Cells.Item(4, 6).Value = -11.473420649752
$ws.Cells.Item(4, 7).Value = -8.958975649751999
$ws.Cells.Item(4, 8).Value = -9.888697649751998
$ws.Cells.Item(4, 9).Value = -9.653689649752
$ws.Cells.Item(4, 10).Value = -9.960825649752
$ws.Cells.Item(4, 11).Value = -9.450149649751999
$ws.Cells.Item(5, 2).Value = 8.382980150385469
$ws.Cells.Item(5, 3).Value = 10.17128915038547
$ws.Cells.Item(5, 4).Value = 6.99317015038547
$ws.Cells.Item(5, 5).Value = 6.69784115038547
$ws.Cells.Item(5, 6).Value = 9.21228615038547
$ws.Cells.Item(5, 7).Value = 8.28256415038547
$ws.Cells.Item(5, 8).Value = 8.517572150385469
$ws.Cells.Item(5, 9).Value = 8.210436150385469
$ws.Cells.Item(5, 10).Value = 8.721112150385469
$ws.Cells.Item(5, 11).Value = 8.433745150385469
$ws.Cells.Item(6, 2).Value = 2.03161760340632
$ws.Cells.Item(6, 3).Value = -1.146501396593679
$ws.Cells.Item(6, 4).Value = -1.441830396593679
$ws.Cells.Item(6, 5).Value = 1.072614603406321
$ws.Cells.Item(6, 6).Value = 0.1428926034063205
$ws.Cells.Item(6, 7).Value = 0.3779006034063205
$ws.Cells.Item(6, 8).Value = 0.07076460340632051
$ws.Cells.Item(6, 9).Value = 0.5814406034063205
$ws.Cells.Item(6, 10).Value = 0.2940736034063205
$ws.Cells.Item(6, 11).Value = 0.3220726034063205
$ws.Cells.Item(7, 2).Value = -3.294071868393072
$ws.Cells.Item(7, 3).Value = -3.589400868393072
$ws.Cells.Item(7, 4).Value = -1.074955868393072
$ws.Cells.Item(7, 5).Value = -2.004677868393072
$ws.Cells.Item(7, 6).Value = -1.769669868393072
$ws.Cells.Item(7, 7).Value = -2.076805868393072
$ws.Cells.Item(7, 8).Value = -1.566129868393072
$ws.Cells.Item(7, 9).Value = -1.853496868393072
$ws.Cells.Item(7, 10).Value = -1.825497868393072
$ws.Cells.Item(7, 11).Value = -1.732849868393072
$ws.Cells.Item(8, 2).Value = -0.1493708818697419
$ws.Cells.Item(8, 3).Value = 2.365074118130258
$ws.Cells.Item(8, 4).Value = 1.435352118130258
$ws.Cells.Item(8, 5).Value = 1.670360118130258
$ws.Cells.Item(8, 6).Value = 1.363224118130258
$ws.Cells.Item(8, 7).Value = 1.873900118130258
$ws.Cells.Item(8, 8).Value = 1.586533118130258
$ws.Cells.Item(8, 9).Value = 1.614532118130258
$ws.Cells.Item(8, 10).Value = 1.707180118130258
$ws.Cells.Item(8, 11).Value = 1.920156118130258
$ws.Cells.Item(9, 2).Value = 2.426248296544459
$ws.Cells.Item(9, 3).Value = 1.496526296544459
$ws.Cells.Item(9, 4).Value = 1.731534296544459
$ws.Cells.Item(9, 5).Value = 1.424398296544459
$ws.Cells.Item(9, 6).Value = 1.935074296544459
$ws.Cells.Item(9, 7).Value = 1.647707296544459
$ws.Cells.Item(9, 8).Value = 1.675706296544459
$ws.Cells.Item(9, 9).Value = 1.768354296544459
$ws.Cells.Item(9, 10).Value = 1.981330296544459
$ws.Cells.Item(9, 11).Value = 1.480670296544459
$ws.Cells.Item(10, 2).Value = -0.5449296406117954
$ws.Cells.Item(10, 3).Value = -0.3099216406117954
$ws.Cells.Item(10, 4).Value = -0.6170576406117954
$ws.Cells.Item(10, 5).Value = -0.1063816406117954
$ws.Cells.Item(10, 6).Value = -0.3937486406117954
$ws.Cells.Item(10, 7).Value = -0.3657496406117954
$ws.Cells.Item(10, 8).Value = -0.2731016406117954
$ws.Cells.Item(10, 9).Value = -0.06012564061179543
$ws.Cells.Item(10, 10).Value = -0.5607856406117955
$ws.Cells.Item(10, 11).Value = -0.2710286406117954
$ws.Cells.Item(11, 2).Value = 0.2726728964218456
$ws.Cells.Item(11, 3).Value = -0.03446310357815441
$ws.Cells.Item(11, 4).Value = 0.4762128964218456
$ws.Cells.Item(11, 5).Value = 0.1888458964218456
$ws.Cells.Item(11, 6).Value = 0.2168448964218456
$ws.Cells.Item(11, 7).Value = 0.3094928964218456
$ws.Cells.Item(11, 8).Value = 0.5224688964218456
$ws.Cells.Item(11, 9).Value = 0.02180889642184558
$ws.Cells.Item(11, 10).Value = 0.3115658964218456
$ws.Cells.Item(11, 11).Value = 0.05367389642184558
$ws.Cells.Item(12, 2).Value = -0.1301381443875122
$ws.Cells.Item(12, 3).Value = 0.3805378556124878
$ws.Cells.Item(12, 4).Value = 0.09317085561248779
$ws.Cells.Item(12, 5).Value = 0.1211698556124878
$ws.Cells.Item(12, 6).Value = 0.2138178556124878
$ws.Cells.Item(12, 7).Value = 0.4267938556124878
$ws.Cells.Item(12, 8).Value = -0.07386614438751221
$ws.Cells.Item(12, 9).Value = 0.2158908556124878
$ws.Cells.Item(12, 10).Value = -0.04200114438751221
$ws.Cells.Item(12, 11).Value = 0.3403798556124878
$ws.Cells.Item(13, 2).Value = 0.6409568926112106
$ws.Cells.Item(13, 3).Value = 0.3535898926112106
$ws.Cells.Item(13, 4).Value = 0.3815888926112106
$ws.Cells.Item(13, 5).Value = 0.4742368926112106
$ws.Cells.Item(13, 6).Value = 0.6872128926112107
$ws.Cells.Item(13, 7).Value = 0.1865528926112106
$ws.Cells.Item(13, 8).Value = 0.4763098926112106
$ws.Cells.Item(13, 9).Value = 0.2184178926112106
$ws.Cells.Item(13, 10).Value = 0.6007988926112107
$ws.Cells.Item(13, 11).Value = 0.1895217986112106
$ws.Cells.Item(14, 2).Value = -0.4818651035472806
$ws.Cells.Item(14, 3).Value = -0.4538661035472806
$ws.Cells.Item(14, 4).Value = -0.3612181035472806
$ws.Cells.Item(14, 5).Value = -0.1482421035472806
$ws.Cells.Item(14, 6).Value = -0.6489021035472806
$ws.Cells.Item(14, 7).Value = -0.3591451035472806
$ws.Cells.Item(14, 8).Value = -0.6170371035472806
$ws.Cells.Item(14, 9).Value = -0.2346561035472806
$ws.Cells.Item(14, 10).Value = -0.6459331975472806
$ws.Cells.Item(14, 11).Value = -0.3589771035472806
$ws.Cells.Item(15, 2).Value = -0.3537867436446591
$ws.Cells.Item(15, 3).Value = -0.2611387436446591
$ws.Cells.Item(15, 4).Value = -0.04816274364465911
$ws.Cells.Item(15, 5).Value = -0.5488227436446591
$ws.Cells.Item(15, 6).Value = -0.2590657436446591
$ws.Cells.Item(15, 7).Value = -0.5169577436446591
$ws.Cells.Item(15, 8).Value = -0.1345767436446591
$ws.Cells.Item(15, 9).Value = -0.5458538376446591
$ws.Cells.Item(15, 10).Value = -0.2588977436446591
$ws.Cells.Item(15, 11).ClearContents()
$ws.Cells.Item(16, 2).Value = 0.09264776243503714
$ws.Cells.Item(16, 3).Value = 0.3056237624350371
$ws.Cells.Item(16, 4).Value = -0.1950362375649629
$ws.Cells.Item(16, 5).Value = 0.09472076243503715
$ws.Cells.Item(16, 6).Value = -0.1631712375649629
$ws.Cells.Item(16, 7).Value = 0.2192097624350371
$ws.Cells.Item(16, 8).Value = -0.1920673315649629
$ws.Cells.Item(16, 9).Value = 0.09488876243503713
$ws.Cells.Item(16, 10).ClearContents()
$ws.Cells.Item(17, 2).Value = 0.1574463720025918
$ws.Cells.Item(17, 3).Value = -0.3432136279974082
$ws.Cells.Item(17, 4).Value = -0.05345662799740819
$ws.Cells.Item(17, 5).Value = -0.3113486279974082
$ws.Cells.Item(17, 6).Value = 0.0710323720025918
$ws.Cells.Item(17, 7).Value = -0.3402447219974082
$ws.Cells.Item(17, 8).Value = -0.0532886279974082
$ws.Cells.Item(17, 9).ClearContents()
$ws.Cells.Item(18, 2).Value = -0.5006596170015631
$ws.Cells.Item(18, 3).Value = -0.2109026170015632
$ws.Cells.Item(18, 4).Value = -0.4687946170015632
$ws.Cells.Item(18, 5).Value = -0.08641361700156319
$ws.Cells.Item(18, 6).Value = -0.4976907110015631
$ws.Cells.Item(18, 7).Value = -0.2107346170015632
$ws.Cells.Item(18, 8).ClearContents()
$ws.Cells.Item(19, 2).Value = 0.2897568395245076
$ws.Cells.Item(19, 3).Value = 0.03186483952450761
$ws.Cells.Item(19, 4).Value = 0.4142458395245076
$ws.Cells.Item(19, 5).Value = 0.002968745524507627
$ws.Cells.Item(19, 6).Value = 0.2899248395245076
$ws.Cells.Item(19, 7).ClearContents()
$ws.Cells.Item(20, 2).Value = -0.2484930410109615
$ws.Cells.Item(20, 3).Value = 0.1338879589890384
$ws.Cells.Item(20, 4).Value = -0.2773891350109615
$ws.Cells.Item(20, 5).Value = 0.009566958989038449
$ws.Cells.Item(20, 6).ClearContents()
$ws.Cells.Item(21, 2).Value = 0.3058628168340501
$ws.Cells.Item(21, 3).Value = -0.1054142771659499
$ws.Cells.Item(21, 4).Value = 0.1815418168340501
$ws.Cells.Item(21, 5).ClearContents()
$ws.Cells.Item(22, 2).Value = -0.4112768169122814
$ws.Cells.Item(22, 3).Value = -0.1243207229122814
$ws.Cells.Item(22, 4).ClearContents()
$ws.Cells.Item(23, 2).Value = 0.4880092297750048
$ws.Cells.Item(23, 3).ClearContents()
$ws.Cells.Item(24, 2).ClearContents()

Write-Output "done"